$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at M (pushes the old "obs"/"normalization" columns
# right, from M/N to N/O) and fill it in with the new "syst_c" metric.
$ws.Columns("M").Insert()

$ws.Range("M1").Value = "syst_c"
$ws.Range("M2:M29").Formula = "=0.06*J2"
